$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pm = [char]0x00B1

# Set column B values first (rows 2-8), then column C values (rows 2-8),
# so the shared-strings table is rebuilt in the same order as the source.
$ws.Range("B2").Value = "+3789.26 $pm 21.87"
$ws.Range("B3").Value = "+3789.26 $pm 21.87"
$ws.Range("B4").Value = "-3812.70 $pm 21.36"
$ws.Range("B5").Value = "-1898.77 $pm 8.88"
$ws.Range("B6").Value = "-1913.94 $pm 14.32"
$ws.Range("B7").Value = "+41.83 $pm 0.24"
$ws.Range("B8").Value = "+4.00 $pm 0.03"

$ws.Range("C2").Value = "+5510.06 $pm 14.81"
$ws.Range("C3").Value = "+5510.06 $pm 14.81"
$ws.Range("C4").Value = "-5559.09 $pm 14.96"
$ws.Range("C5").Value = "-2327.11 $pm 5.73"
$ws.Range("C6").Value = "-3231.97 $pm 9.45"
$ws.Range("C7").Value = "+9.66 $pm 0.02"
$ws.Range("C8").Value = "+0.63 $pm 0.00"
